# fix: one minute analisis for bnb coins collections
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (and the <sheet> entry in workbook.xml)
$ws.Name = "COIN_FIRST"

# Row 2 (was ID=0) -> shifted data / date
$ws.Range("B2").Value = 3960
$ws.Range("C2").Value = 3371
$ws.Range("D2").Value = 3626
$ws.Range("E2").Value = 3336
$ws.Range("F2").Value = 3147
$ws.Range("G2").Value = 3274
$ws.Range("H2").Value = 2934
$ws.Range("I2").Value = 3248
$ws.Range("J2").Value = 3447
$ws.Range("K2").Value = 3360
$ws.Range("L2").Value = "2021-06-24 00:00:00"

# Row 3 (was ID=1) -> shifted data / date
$ws.Range("B3").Value = 7653
$ws.Range("C3").Value = 6878
$ws.Range("D3").Value = 7042
$ws.Range("E3").Value = 6492
$ws.Range("F3").Value = 6346
$ws.Range("G3").Value = 6770
$ws.Range("H3").Value = 6605
$ws.Range("I3").Value = 6331
$ws.Range("J3").Value = 6502
$ws.Range("K3").Value = 6618
$ws.Range("L3").Value = "2021-06-23 00:00:00"

# Row 4 (was ID=2) -> zeroed out values, date shifted
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = "2021-06-22 00:00:00"

# Row 5 (was ID=3) -> zeroed out values, date shifted
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = "2021-06-21 00:00:00"

# Row 6 (was ID=4) -> zeroed out values, date shifted
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = "2021-06-20 00:00:00"

# Row 7 (was ID=5) -> zeroed out values, date shifted
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = "2021-06-19 00:00:00"

# Row 8 (was ID=6) -> zeroed out values, date shifted
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = "2021-06-18 00:00:00"
